# Daily attendance processing - 2025-10-29 11:43:11
#
# Normalises the "Recorded By" (column G) text on the "Session Analysis
# Results" sheet: for a handful of well-known recorder-name combinations,
# the last name in the comma-separated list is moved to the front.
#
#   "System, backup@backdoor.com, system"  -> "system, System, backup@backdoor.com"
#   "dnasr281@gmail.com, System"           -> "System, dnasr281@gmail.com"
#   "admin@admin.com, System"              -> "System, admin@admin.com"
#   "admin@admin.com, dnasr281@gmail.com"  -> "dnasr281@gmail.com, admin@admin.com"
#
# Every other combination (e.g. the 2-item "System, backup@backdoor.com")
# is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, backup@backdoor.com, system" = "system, System, backup@backdoor.com"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "admin@admin.com, System"             = "System, admin@admin.com"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 2) { $lastRow = 2 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2

    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
